$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("step_change")
$ws2 = $wb.Worksheets.Item("notes")

# --- step_change sheet: insert a new "old" column between the existing
# year/EI(g_wh) columns, shifting the former B (EI data) header+values
# into C and populating the new B column with the "old" header text.

# Header row: B1 becomes "old", C1 becomes the former "EI(g_wh)" header.
$ws.Range("B1").Value = "old"
$ws.Range("C1").Value = "EI(g_wh)"

# New column C data values (rows 2-31).
$ws.Range("C2").Value = 0.62422682874642976
$ws.Range("C3").Value = 0.61939975050351725
$ws.Range("C4").Value = 0.61540940649038212
$ws.Range("C5").Value = 0.59871070537004722
$ws.Range("C6").Value = 0.58387019051458278
$ws.Range("C7").Value = 0.55465722631191772
$ws.Range("C8").Value = 0.57214187305195952
$ws.Range("C9").Value = 0.51823856546817315
$ws.Range("C10").Value = 0.4997748244352076
$ws.Range("C11").Value = 0.50539883731315405
$ws.Range("C12").Value = 0.46637253828755709
$ws.Range("C13").Value = 0.41443359671112523
$ws.Range("C14").Value = 0.41600911249874645
$ws.Range("C15").Value = 0.40469834416909695
$ws.Range("C16").Value = 0.30850420341150847
$ws.Range("C17").Value = 0.27989452129631548
$ws.Range("C18").Value = 0.25264010087669825
$ws.Range("C19").Value = 0.24722333414421613
$ws.Range("C20").Value = 0.24841309236874204
$ws.Range("C21").Value = 0.24474313158762934
$ws.Range("C22").Value = 0.23062618419074904
$ws.Range("C23").Value = 0.20500105261399915
$ws.Range("C24").Value = 0.17937592103724925
$ws.Range("C25").Value = 0.15375078946049936
$ws.Range("C26").Value = 0.12812565788374947
$ws.Range("C27").Value = 0.10250052630699957
$ws.Range("C28").Value = 0.07687539473024968
$ws.Range("C29").Value = 0.051250263153499787
$ws.Range("C30").Value = 0.025625131576749893
$ws.Range("C31").Value = 0

# Rows 23-31 in column B/A use the shaded "Per cent" style; copy that
# formatting onto the matching C cells (keeps the values just written).
$ws.Range("B23:B31").Copy()
$ws.Range("C23:C31").PasteSpecial(-4122)

# Row 32 is new: a single styled-but-empty cell in column C, matching the
# shaded style used by the rows above it.
$ws.Range("B31").Copy()
$ws.Range("C32").PasteSpecial(-4122)
$ws.Range("C32").ClearContents()

$excel.CutCopyMode = $false

# Selection moves to C3 on the step_change sheet.
$ws.Range("C3").Select()

# --- notes sheet: scroll the view so row 5 is the top-left visible cell.
$ws2.Activate()
$ws2.Application.ActiveWindow.ScrollRow = 5
$ws.Activate()
